# Adds the two newest quarters (2024Q2 / 2024Q3) of churn data to the
# bottom of the table on Sheet1, rows 9 and 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 : 2024Q2 -------------------------------------------------
$ws.Cells.Item(9, 1).Value = "'15683481"
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(9, 1).PasteSpecial(-4122)          # xlPasteFormats - keep default (unstyled) look

$ws.Cells.Item(9, 2).Value = 2024
$ws.Cells.Item(9, 3).Value = 175476
$ws.Cells.Item(9, 4).Value = "BPO Løn & HR"
$ws.Cells.Item(9, 5).Value = 45428
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(9, 5).PasteSpecial(-4122)          # xlPasteFormats - reuse the date-style (s="2")
$ws.Cells.Item(9, 8).Value = "2024Q2"
$ws.Cells.Item(9, 9).Value = "160000-180000"

# --- Row 10 : 2024Q3 -------------------------------------------------
$ws.Cells.Item(10, 1).Value = "'25638670"
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(10, 1).PasteSpecial(-4122)

$ws.Cells.Item(10, 2).Value = 2024
$ws.Cells.Item(10, 3).Value = 177066
$ws.Cells.Item(10, 4).Value = "BPO Løn og refusion"
$ws.Cells.Item(10, 5).Value = 45531
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(10, 5).PasteSpecial(-4122)         # xlPasteFormats - reuse the date-style (s="2")
$ws.Cells.Item(10, 8).Value = "2024Q3"
$ws.Cells.Item(10, 9).Value = "160000-180000"

$excel.CutCopyMode = $false
